# Insert a new weekly price-report row for "Ají" (Inferno variety) at
# row 205 of the single data sheet, shifting the existing rows 205-223
# down to 206-224 (dimension grows from A1:R223 to A1:R224).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push rows 205..223 down by one to make room for the new record.
$ws.Rows(205).Insert()

# Populate the newly inserted row 205 with the new market record.
$ws.Range("A205").Value = 7
$ws.Range("B205").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C205").Value = "Ñuble"
$ws.Range("D205").Value = 45142
$ws.Range("E205").Value = 16
$ws.Range("F205").Value = 100112021
$ws.Range("G205").Value = "Ají"
$ws.Range("H205").Value = "Inferno"
$ws.Range("I205").Value = "Primera"
$ws.Range("J205").Value = 30
$ws.Range("K205").Value = 14000
$ws.Range("L205").Value = 14000
$ws.Range("M205").Value = 14000
$ws.Range("N205").Value = "$/caja 10 kilos"
$ws.Range("O205").Value = "Región de Arica y Parinacota"
$ws.Range("P205").Value = 1400
$ws.Range("Q205").Value = 10
$ws.Range("R205").Value = "Hortaliza"
